$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.023556
$ws.Range("H2").Value = 0.07066799999999999
$ws.Range("I2").Value = 0.08088674619362546
$ws.Range("J2").Value = 0.08088674619362549
$ws.Range("Q2").Value = 0.0008080493199999998
$ws.Range("R2").Value = 0.00727244388
$ws.Range("S2").Value = 0.08088674619362546
$ws.Range("T2").Value = 0.08088674619362549

# Row 3
$ws.Range("I3").Value = 0.8589438069010353
$ws.Range("J3").Value = 0.8589438069010354
$ws.Range("R3").Value = 0.07722675129999999
$ws.Range("S3").Value = 0.8589438069010353
$ws.Range("T3").Value = 0.8589438069010354

# Row 4
$ws.Range("G4").Value = 0.01752266666666667
$ws.Range("H4").Value = 0.052568
$ws.Range("I4").Value = 0.0601694469053391
$ws.Range("J4").Value = 0.06016944690533912
$ws.Range("Q4").Value = 0.0006010858755555555
$ws.Range("R4").Value = 0.00540977288
$ws.Range("S4").Value = 0.0601694469053391
$ws.Range("T4").Value = 0.06016944690533912
